# Append one new data row (row 19) to Sheet1, mirroring the existing
# date/OHLCV/headline columns: date, close_price, high_price, low_price,
# open_price, volume, daily_return, label, daily_headlines, CleanText.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

# date (serial 45859 == 2025-07-21) - match formatting used by the other
# rows in column A (same date/time number format).
$ws.Cells.Item($row, 1).Value = 45859
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# close_price, high_price, low_price, open_price, volume
$ws.Cells.Item($row, 2).Value = 6305.60009765625
$ws.Cells.Item($row, 3).Value = 6336.080078125
$ws.Cells.Item($row, 4).Value = 6303.7900390625
$ws.Cells.Item($row, 5).Value = 6304.740234375
$ws.Cells.Item($row, 6).Value = 5010840000

# daily_return, label
$ws.Cells.Item($row, 7).Value = 0.0013991348828683
$ws.Cells.Item($row, 8).Value = 1

# daily_headlines (empty for this date), CleanText
$ws.Cells.Item($row, 9).Value = ""
$ws.Cells.Item($row, 10).Value = "nan"
